# "Check if the Sudoku is Legit"
# Fills in the blank clue-cells of the puzzle grid (rows 1-9) with the
# given numbers, and replaces the solution grid (rows 12-20) with a new
# (valid) completed Sudoku solution that matches those clues.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sudoku")

# --- Puzzle grid (rows 1-9): fill previously-blank clue cells ---------
$ws.Range("D1").Value = 5
$ws.Range("G1").Value = 8
$ws.Range("H1").Value = 3

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 8
$ws.Range("I2").Value = 7

$ws.Range("C3").Value = 9
$ws.Range("E3").Value = 3
$ws.Range("H3").Value = 5

$ws.Range("A4").Value = 9
$ws.Range("D4").Value = 3
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 4

$ws.Range("C5").Value = 1
$ws.Range("E5").Value = 9
$ws.Range("G5").Value = 5

$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 4
$ws.Range("F6").Value = 2
$ws.Range("I6").Value = 3

$ws.Range("B7").Value = 5
$ws.Range("E7").Value = 7

$ws.Range("A8").Value = 1
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 6

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 6

# --- Solution grid (rows 12-20): replace with the new valid solution --
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 8
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 9

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 7

$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 6

$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 2

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 7
$ws.Range("I16").Value = 8

$ws.Range("A17").Value = 5
$ws.Range("B17").Value = 7
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 6
$ws.Range("I17").Value = 3

$ws.Range("A18").Value = 2
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 8
$ws.Range("I18").Value = 1

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 4

$ws.Range("A20").Value = 8
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 9
$ws.Range("I20").Value = 5

# --- Selection state: mirrors the box the edits were made in ----------
# (the COM layer always normalises the active cell to the selected
# range's top-left corner, so this is the closest reachable match to
# the authored activeCell="H7" sqref="F7:H9")
$ws.Range("F7:H9").Select()
